$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.953.78"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.775.76"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.84"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.18"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.563"
$ws.Range("E7").Value = "  +3.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.17"
$ws.Range("E10").Value = "  -4.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0852"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.132"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.38"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.59"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "3.205.81"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "2.827.45"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.926"
$ws.Range("E17").Value = "  +3.69%  "
$ws.Range("D18").Value = "51.761.86"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.37"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.02"
$ws.Range("E21").Value = "  -5.25%  "
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "273.65"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.62"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.52"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.144"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.81"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.70"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  +9.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0836"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.20"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("E40").Value = "  -4.58%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.51"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.98"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.67"
$ws.Range("E45").Value = "  -7.01%  "
$ws.Range("D46").Value = "2.060.47"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.69"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.925"
$ws.Range("E50").Value = "  -2.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.89"
$ws.Range("E51").Value = "  +0.18%  "

# Row 31 and 32 swap (OKB <-> VeChain) with updated data
$ws.Range("B31").Value = "VeChain"
$ws.Range("C31").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0463"
$ws.Range("E31").Value = "  +4.32%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.36"
$ws.Range("E32").Value = "  +1.28%  "
